# The deck currently ships two theme parts:
#   ppt/theme/theme1.xml  -> "Office Theme" (standard Office colours) - linked from the notes master
#   ppt/theme/theme2.xml  -> "Integral"     (custom green/gold palette) - linked from the slide master
#
# The authored change swaps the two themes' contents, so the slide master
# (and therefore every slide) ends up using the plain "Office Theme" colour
# palette instead of "Integral".
#
# Drive this through the theme colour scheme that every slide shares with
# the slide master: rewriting it propagates straight into the master's
# theme part (ppt/theme/theme2.xml) for the whole deck at once.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# ppThemeColorDark1 .. ppThemeColorFollowedHyperlink (indices 1-12), set to
# the stock "Office Theme" palette (the colours theme1.xml carried before
# the swap).
$tcs.Item(1).RGB  = 0           # Dark 1     -> 000000
$tcs.Item(2).RGB  = 16777215    # Light 1    -> FFFFFF
$tcs.Item(3).RGB  = 6968388     # Dark 2     -> 44546A
$tcs.Item(4).RGB  = 15132391    # Light 2    -> E7E6E6
$tcs.Item(5).RGB  = 13998939    # Accent 1   -> 5B9BD5
$tcs.Item(6).RGB  = 3243501     # Accent 2   -> ED7D31
$tcs.Item(7).RGB  = 10855845    # Accent 3   -> A5A5A5
$tcs.Item(8).RGB  = 49407       # Accent 4   -> FFC000
$tcs.Item(9).RGB  = 12874308    # Accent 5   -> 4472C4
$tcs.Item(10).RGB = 4697456     # Accent 6   -> 70AD47
$tcs.Item(11).RGB = 12673797    # Hyperlink  -> 0563C1
$tcs.Item(12).RGB = 7491477     # Followed Hyperlink -> 954F72
